$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns keep their text formatting so numeric-looking
# strings (prices like "1.00", percentages) are not coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "87.724.76"
$ws.Range("E2").Value = "  +8.06%  "
$ws.Range("D3").Value = "3.330.11"
$ws.Range("E3").Value = "  +4.57%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "219.15"
$ws.Range("E5").Value = "  +4.54%  "
$ws.Range("D6").Value = "652.15"
$ws.Range("E6").Value = "  +2.57%  "
$ws.Range("D7").Value = "0.356"
$ws.Range("E7").Value = "  +22.81%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("D10").Value = "3.329.00"
$ws.Range("E10").Value = "  +4.72%  "
$ws.Range("D11").Value = "0.588"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").Value = "0.0000270"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "35.44"
$ws.Range("E14").Value = "  +10.20%  "
$ws.Range("D15").Value = "3.938.00"
$ws.Range("E15").Value = "  +4.57%  "
$ws.Range("D16").Value = "5.49"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").Value = "87.566.39"
$ws.Range("E17").Value = "  +7.91%  "
$ws.Range("D18").Value = "3.321.90"
$ws.Range("E18").Value = "  +4.70%  "
$ws.Range("D19").Value = "14.78"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").Value = "9.88"
$ws.Range("E20").Value = "  +6.57%  "
$ws.Range("D21").Value = "3.15"
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("D22").Value = "457.17"
$ws.Range("E22").Value = "  +3.59%  "
$ws.Range("D23").Value = "5.56"
$ws.Range("E23").Value = "  +6.05%  "
$ws.Range("D24").Value = "5.62"
$ws.Range("E24").Value = "  +10.32%  "
$ws.Range("D25").Value = "12.73"
$ws.Range("E25").Value = "  +12.43%  "
$ws.Range("D26").Value = "3.488.40"
$ws.Range("E26").Value = "  +4.14%  "
$ws.Range("D27").Value = "79.09"
$ws.Range("E27").Value = "  +2.71%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("B29").Value = "Cronos"
$ws.Range("C29").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D29").Value = "0.196"
$ws.Range("E29").Value = "  +37.21%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0000127"
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "9.46"
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "607.17"
$ws.Range("E32").Value = "  +5.77%  "
$ws.Range("E33").Value = "  +6.21%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("D35").Value = "2.11"
$ws.Range("E35").Value = "  +3.02%  "
$ws.Range("D36").Value = "7.16"
$ws.Range("E36").Value = "  +19.74%  "
$ws.Range("E37").Value = "  -3.84%  "
$ws.Range("D38").Value = "23.38"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").Value = "2.15"
$ws.Range("E39").Value = "  +4.14%  "
$ws.Range("D40").Value = "0.422"
$ws.Range("E40").Value = "  +1.70%  "
$ws.Range("D41").Value = "21.84"
$ws.Range("E41").Value = "  +5.01%  "
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "3.04"
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("D44").Value = "159.58"
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "191.94"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("D47").Value = "1.42"
$ws.Range("E47").Value = "  +5.34%  "
$ws.Range("D48").Value = "46.25"
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D49").Value = "4.49"
$ws.Range("E49").Value = "  +4.20%  "
$ws.Range("D50").Value = "0.786"
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("D51").Value = "0.664"
$ws.Range("E51").Value = "  +2.50%  "
